# Update the "Förändrad" (Changed) date column (C) for all data rows (2-36)
# from serial date 45686 to 45687 (i.e. +1 day), matching the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 36; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45686) {
        $cell.Value2 = 45687
    }
}
